# New crime data collected — weekly refresh of the 111th Precinct CompStat
# report: bump the report volume/week-ending text, and roll the Week-to-Date /
# 28-Day / Year-to-Date figures (and their derived % changes) forward by one
# week for Robbery, Fel. Assault, Burglary, Gr. Larceny, G.L.A., TOTAL,
# Petit Larceny, Misd. Assault, Other Sex Crimes and Hate Crimes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number + the reporting week's date range ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Row 16: Robbery ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 30
$ws.Range("L16").Value = -7.142857142857
$ws.Range("M16").Value = -23.529411764705
$ws.Range("N16").Value = -65.789473684210

# --- Row 17: Fel. Assault ---
$ws.Range("D17").Value = 4
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 22.222222222222
$ws.Range("L17").Value = -15.384615384615
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -35.294117647058

# --- Row 18: Burglary ---
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = -22.222222222222
$ws.Range("F18").Value = 21
$ws.Range("H18").Value = -34.375
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -39.622641509434
$ws.Range("L18").Value = -28.888888888888
$ws.Range("M18").Value = -15.789473684210
$ws.Range("N18").Value = -77.777777777777

# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -45.205479452054
$ws.Range("I19").Value = 59
$ws.Range("J19").Value = 113
$ws.Range("K19").Value = -47.787610619469
$ws.Range("L19").Value = -14.492753623188
$ws.Range("M19").Value = 3.508771929824
$ws.Range("N19").Value = -10.606060606060

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 36
$ws.Range("J20").Value = 26
$ws.Range("K20").Value = 38.461538461538
$ws.Range("L20").Value = 125
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = -91.529411764705

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = -35.074626865671
$ws.Range("I21").Value = 152
$ws.Range("J21").Value = 211
$ws.Range("K21").Value = -27.962085308056
$ws.Range("L21").Value = -3.184713375796
$ws.Range("M21").Value = 10.144927536231
$ws.Range("N21").Value = -78.002894356005

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 44
$ws.Range("H24").Value = -12
$ws.Range("I24").Value = 89
$ws.Range("J24").Value = 77
$ws.Range("K24").Value = 15.584415584415
$ws.Range("L24").Value = -21.929824561403
$ws.Range("M24").Value = 3.488372093023

# --- Row 25: Misd. Assault ---
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 25
$ws.Range("J25").Value = 26
$ws.Range("K25").Value = -3.846153846153
$ws.Range("L25").Value = -28.571428571428
$ws.Range("M25").Value = 19.047619047619

# --- Row 27: Other Sex Crimes ---
# D27/E27/L27 flip from the "N/A" / "***.*" placeholder text to real numbers,
# so give them the same numeric formats the rest of the table uses.
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J27").Value = 2
$ws.Range("L27").Value = -100
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Row 30: Hate Crimes ---
$ws.Range("L30").Value = 0
$ws.Range("L30").NumberFormat = "#,##0.0;""-""#,##0.0"
